$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 439, shifting rows 439:481 down to 440:482.
$ws.Rows.Item(439).Insert()

# Populate the newly inserted row 439 with the new price record.
$ws.Cells.Item(439, 1).Value = 7
$ws.Cells.Item(439, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(439, 3).Value = "Ñuble"
$ws.Cells.Item(439, 4).Value = 45132
$ws.Cells.Item(439, 5).Value = 16
$ws.Cells.Item(439, 6).Value = 100112003
$ws.Cells.Item(439, 7).Value = "Ajo"
$ws.Cells.Item(439, 8).Value = "Chino"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 50
$ws.Cells.Item(439, 11).Value = 21000
$ws.Cells.Item(439, 12).Value = 21000
$ws.Cells.Item(439, 13).Value = 21000
$ws.Cells.Item(439, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(439, 15).Value = "China"
$ws.Cells.Item(439, 16).Value = 2100
$ws.Cells.Item(439, 17).Value = 10
$ws.Cells.Item(439, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date/time number format as its
# neighbours (style index "2" -> numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(439, 4).NumberFormat = $ws.Cells.Item(440, 4).NumberFormat
